$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as text (inline strings)
# in the workbook. Writing a plain numeric-looking string to a cell causes Excel
# to convert it to a real number (losing formatting such as trailing zeros, e.g.
# "0.09300" -> 0.093). Temporarily formatting the cell as Text ("@") while we set
# the value keeps it as text, and resetting the style back to Normal afterwards
# avoids leaving a residual custom style assigned to the cell.
function Set-TextValue($address, $value) {
    $range = $ws.Range($address)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue "D2" "247.44"
Set-TextValue "D3" "21.71"
Set-TextValue "D4" "5.344"
Set-TextValue "D5" "0.05628"
Set-TextValue "D6" "3.431"
Set-TextValue "D7" "6.372"
Set-TextValue "D8" "0.8125"
Set-TextValue "D9" "0.9359"
Set-TextValue "D10" "0.1425"
Set-TextValue "D11" "0.07537"
Set-TextValue "D12" "0.03217"
Set-TextValue "D13" "0.03087"
Set-TextValue "D14" "0.09300"
Set-TextValue "D15" "3.595"
Set-TextValue "D16" "0.001587"
Set-TextValue "D17" "0.04718"
Set-TextValue "D18" "0.0005783"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue "D19" "0.006293"
Set-TextValue "D21" "0.001032"
Set-TextValue "D22" "0.0001500"
Set-TextValue "D23" "3.762"
Set-TextValue "D25" "0.3301"
Set-TextValue "D40" "0.03953"
Set-TextValue "D41" "0.1065"
Set-TextValue "D42" "0.003020"
Set-TextValue "D43" "0.002926"
Set-TextValue "D44" "0.008799"
Set-TextValue "D45" "0.00005584"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.0005503"
$ws.Range("E47").Value = "46ACDXExchangeACXTWorstin24h"
Set-TextValue "D48" "0.7804"
Set-TextValue "D49" "0.1769"
